{"js": "// The document contains a single paragraph run with the text\n// \", ch\u00fang t\u00f4i g\u1ed3m:\" which must become \", g\u1ed3m c\u00f3:\" (i.e. the word\n// \"ch\u00fang t\u00f4i\" is dropped and \"c\u00f3\" is inserted before the trailing colon).\nconst body = context.document.body;\nconst results = body.search(\", ch\u00fang t\u00f4i g\u1ed3m:\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text \", ch\u00fang t\u00f4i g\u1ed3m:\" not found.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\", g\u1ed3m c\u00f3:\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The paragraph run \", ch\u00fang t\u00f4i g\u1ed3m:\" must become \", g\u1ed3m c\u00f3:\" (the\n# word \"ch\u00fang t\u00f4i\" is removed and \"c\u00f3\" is inserted before the colon).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$findText = \", ch\u00fang t\u00f4i g\u1ed3m:\"\n$replaceText = \", g\u1ed3m c\u00f3:\"\n\n$find.Text = $findText\n$find.Replacement.Text = $replaceText\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\nif (-not $result) {\n    throw 'Target text \", chung toi gom:\" not found.'\n}\n"}
